# Generate Report for handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) values for the first data row
# (the 23683442-... file) on both the zh-cn and de-de sheets, reflecting a
# newer handoff/handback cycle.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-02-17 02:59:45"
$wsZhCn.Range("G2").Value = "2016-02-17 03:00:36"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-02-17 02:59:55"
$wsDeDe.Range("G2").Value = "2016-02-17 03:00:57"
